# Insert a new row at position 21 (shifts old rows 21-50 down to 22-51)
# and populate the new row 21 with the "La Notte degli archivi" event.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Cells.Item(21, 1).Value  = "Altri eventi,Visite guidate"
$ws.Cells.Item(21, 2).Value  = "Modena"
$ws.Cells.Item(21, 3).Value  = "vedi nel programma"
$ws.Cells.Item(21, 4).Value  = "2022-06-08T09:03:08+00:00"
$ws.Cells.Item(21, 5).Value  = "Nove archivi modenesi, pubblici e privati, in rete per proporre iniziative e per scoprire le innumerevoli storie che custodiscono"
$ws.Cells.Item(21, 6).Value  = "2022-06-08T09:03:43+00:00"
$ws.Cells.Item(21, 8).Value  = "2022-06-10T09:00:00+00:00"
$ws.Cells.Item(21, 9).Value  = "2022-06-10T10:00:00+00:00"
$ws.Cells.Item(21, 10).Value = "https://www.comune.modena.it/api/novita/eventi/2022/la-notte-degli-archivi/@@images/a4a15f79-4ef8-4285-bb18-6d0019689efe.jpeg"
$ws.Cells.Item(21, 12).Value = "2022-06-08T09:03:43+00:00"
$ws.Cells.Item(21, 13).Value = "Sedi diverse"
$ws.Cells.Item(21, 14).Value = " vedi orari nel programma"
$ws.Cells.Item(21, 19).Value = "La Notte degli archivi"
$ws.Cells.Item(21, 22).Value = $false
$ws.Cells.Item(21, 23).Value = 41123
$ws.Cells.Item(21, 24).Value = "https://www.comune.modena.it/novita/eventi/2022/la-notte-degli-archivi"
$ws.Cells.Item(21, 25).Value = "44,64582"
$ws.Cells.Item(21, 26).Value = "10,92572"
$ws.Cells.Item(21, 27).Value = "POINT (10.92572 44.64582)"
